# Auto-generated edit script: updates Leve profit calculation values across multiple sheets
# per scheduled runner refresh of market price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 398.53
$ws.Range("I15").Value = 398.53
$ws.Range("K15").Value = 1195.59
$ws.Range("M15").Value = -1026.59

$ws.Range("H116").Value = 5258.154
$ws.Range("I116").Value = 2725
$ws.Range("J116").Value = 5718.727
$ws.Range("K116").Value = 2725
$ws.Range("L116").Value = 5718.727
$ws.Range("M116").Value = 717
$ws.Range("N116").Value = -12602.727

$ws.Range("H132").Value = 2756.8857
$ws.Range("I132").Value = 3232.963
$ws.Range("J132").Value = 1150.125
$ws.Range("K132").Value = 9698.889000000001
$ws.Range("L132").Value = 3450.375
$ws.Range("M132").Value = -7168.889000000001
$ws.Range("N132").Value = -8510.375

$ws.Range("H137").Value = 1080.34
$ws.Range("I137").Value = 980.1081
$ws.Range("J137").Value = 1365.6154
$ws.Range("K137").Value = 2940.3243
$ws.Range("L137").Value = 4096.8462
$ws.Range("M137").Value = -390.3243000000002
$ws.Range("N137").Value = -9196.8462

$ws.Range("H138").Value = 19610044
$ws.Range("I138").Value = 31251150
$ws.Range("J138").Value = 3968.842
$ws.Range("K138").Value = 93753450
$ws.Range("L138").Value = 11906.526
$ws.Range("M138").Value = -93748310
$ws.Range("N138").Value = -22186.526

$ws.Range("H141").Value = 1153.7188
$ws.Range("I141").Value = 740.0714
$ws.Range("K141").Value = 2220.2142
$ws.Range("M141").Value = 2959.7858


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3361.3086
$ws.Range("I32").Value = 2907.0898
$ws.Range("J32").Value = 15171
$ws.Range("K32").Value = 2907.0898
$ws.Range("L32").Value = 15171
$ws.Range("M32").Value = -2620.0898
$ws.Range("N32").Value = -15745

$ws.Range("H61").Value = 340665.12
$ws.Range("I61").Value = 400992.47
$ws.Range("J61").Value = 1323.75
$ws.Range("K61").Value = 400992.47
$ws.Range("L61").Value = 1323.75
$ws.Range("M61").Value = -400780.47
$ws.Range("N61").Value = -1747.75

$ws.Range("H63").Value = 2233260.8
$ws.Range("I63").Value = 1203.8462
$ws.Range("K63").Value = 1203.8462
$ws.Range("M63").Value = -517.8462

$ws.Range("H66").Value = 2233260.8
$ws.Range("I66").Value = 1203.8462
$ws.Range("K66").Value = 6019.231
$ws.Range("M66").Value = -2587.231

$ws.Range("H74").Value = 28573646
$ws.Range("I74").Value = 37039164
$ws.Range("K74").Value = 37039164
$ws.Range("M74").Value = -37038290

$ws.Range("H77").Value = 28573646
$ws.Range("I77").Value = 37039164
$ws.Range("K77").Value = 185195820
$ws.Range("M77").Value = -185191452

$ws.Range("H136").Value = 340665.12
$ws.Range("I136").Value = 400992.47
$ws.Range("J136").Value = 1323.75
$ws.Range("K136").Value = 1202977.41
$ws.Range("L136").Value = 3971.25
$ws.Range("M136").Value = -1200427.41
$ws.Range("N136").Value = -9071.25


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 714.8261
$ws.Range("I16").Value = 652.6429000000001
$ws.Range("J16").Value = 811.55554
$ws.Range("K16").Value = 652.6429000000001
$ws.Range("L16").Value = 811.55554
$ws.Range("M16").Value = -365.6429000000001
$ws.Range("N16").Value = -1385.55554

$ws.Range("H31").Value = 5719
$ws.Range("I31").Value = 11706
$ws.Range("J31").Value = 4677.7827
$ws.Range("K31").Value = 11706
$ws.Range("L31").Value = 4677.7827
$ws.Range("M31").Value = -11411
$ws.Range("N31").Value = -5267.7827

$ws.Range("H34").Value = 5719
$ws.Range("I34").Value = 11706
$ws.Range("J34").Value = 4677.7827
$ws.Range("K34").Value = 11706
$ws.Range("L34").Value = 4677.7827
$ws.Range("M34").Value = -11504
$ws.Range("N34").Value = -5081.7827

$ws.Range("H58").Value = 7611.838
$ws.Range("I58").Value = 766.383
$ws.Range("J58").Value = 19528
$ws.Range("K58").Value = 766.383
$ws.Range("L58").Value = 19528
$ws.Range("M58").Value = -563.383
$ws.Range("N58").Value = -19934

$ws.Range("H105").Value = 7813473
$ws.Range("I105").Value = 8929448
$ws.Range("J105").Value = 1650
$ws.Range("K105").Value = 8929448
$ws.Range("L105").Value = 1650
$ws.Range("M105").Value = -8927701
$ws.Range("N105").Value = -5144

$ws.Range("H113").Value = 714.8261
$ws.Range("I113").Value = 652.6429000000001
$ws.Range("J113").Value = 811.55554
$ws.Range("K113").Value = 652.6429000000001
$ws.Range("L113").Value = 811.55554
$ws.Range("M113").Value = 1517.3571
$ws.Range("N113").Value = -5151.55554

$ws.Range("H132").Value = 2036.8431
$ws.Range("I132").Value = 1433.4
$ws.Range("K132").Value = 4300.200000000001
$ws.Range("M132").Value = -1770.200000000001

$ws.Range("H134").Value = 760.3019
$ws.Range("I134").Value = 649.7560999999999
$ws.Range("J134").Value = 1138
$ws.Range("K134").Value = 1949.2683
$ws.Range("L134").Value = 3414
$ws.Range("M134").Value = 585.7317000000003
$ws.Range("N134").Value = -8484

$ws.Range("H136").Value = 7611.838
$ws.Range("I136").Value = 766.383
$ws.Range("J136").Value = 19528
$ws.Range("K136").Value = 2299.149
$ws.Range("L136").Value = 58584
$ws.Range("M136").Value = 250.8509999999997
$ws.Range("N136").Value = -63684

$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1757.4
$ws.Range("I5").Value = 1122.4445
$ws.Range("K5").Value = 3367.3335
$ws.Range("M5").Value = -3255.3335

$ws.Range("H87").Value = 18503.48
$ws.Range("I87").Value = 5999.875
$ws.Range("J87").Value = 24387.53
$ws.Range("K87").Value = 17999.625
$ws.Range("L87").Value = 73162.59
$ws.Range("M87").Value = -16751.625
$ws.Range("N87").Value = -75658.59

$ws.Range("H90").Value = 18503.48
$ws.Range("I90").Value = 5999.875
$ws.Range("J90").Value = 24387.53
$ws.Range("K90").Value = 53998.875
$ws.Range("L90").Value = 219487.77
$ws.Range("M90").Value = -47758.875
$ws.Range("N90").Value = -231967.77

$ws.Range("H120").Value = 12732.728
$ws.Range("I120").Value = 6676.6665
$ws.Range("K120").Value = 20029.9995
$ws.Range("M120").Value = -15191.9995

$ws.Range("H131").Value = 164726.78
$ws.Range("J131").Value = 176226.9
$ws.Range("L131").Value = 528680.7
$ws.Range("N131").Value = -538760.7

$ws.Range("H135").Value = 1757.4
$ws.Range("I135").Value = 1122.4445
$ws.Range("K135").Value = 10102.0005
$ws.Range("M135").Value = -7567.0005


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3481087.8
$ws.Range("I70").Value = 3297.1428
$ws.Range("J70").Value = 5694227
$ws.Range("K70").Value = 3297.1428
$ws.Range("L70").Value = 5694227
$ws.Range("M70").Value = -3027.1428
$ws.Range("N70").Value = -5694767

$ws.Range("H73").Value = 3481087.8
$ws.Range("I73").Value = 3297.1428
$ws.Range("J73").Value = 5694227
$ws.Range("K73").Value = 3297.1428
$ws.Range("L73").Value = 5694227
$ws.Range("M73").Value = -2361.1428
$ws.Range("N73").Value = -5696099

$ws.Range("H122").Value = 121214560
$ws.Range("I122").Value = 47620736
$ws.Range("J122").Value = 250003740
$ws.Range("K122").Value = 142862208
$ws.Range("L122").Value = 750011220
$ws.Range("M122").Value = -142859758
$ws.Range("N122").Value = -750016120

$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2260999.8
$ws.Range("I2").Value = 3501666.8
$ws.Range("J2").Value = 399999.5
$ws.Range("K2").Value = 3501666.8
$ws.Range("L2").Value = 399999.5
$ws.Range("M2").Value = -3501554.8
$ws.Range("N2").Value = -400223.5

$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1370.9697
$ws.Range("I132").Value = 951.875
$ws.Range("J132").Value = 2488.5557
$ws.Range("K132").Value = 2855.625
$ws.Range("L132").Value = 7465.6671
$ws.Range("M132").Value = -325.625
$ws.Range("N132").Value = -12525.6671

